$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 13.03
$ws.Range("E6").Value = 12.642
$ws.Range("D7").Value = -7.284000000000001
$ws.Range("E7").Value = 12.648
$ws.Range("B8").Value = 6.156000000000001
$ws.Range("E8").Value = 12.147
$ws.Range("E9").Value = 12.547
$ws.Range("B10").Value = 7.491000000000001
$ws.Range("E10").Value = 12.143
$ws.Range("B12").Value = 6.444
$ws.Range("E12").Value = 13.055
$ws.Range("C13").Value = -12.217
$ws.Range("B18").Value = 6.258999999999999
$ws.Range("D20").Value = -8.222
$ws.Range("B25").Value = 6.407999999999999
